# Apply updated crypto price/volume data per diff.
# Values are prefixed with a leading apostrophe and the cell style is reset
# to "Normal" afterwards so Excel stores them as plain text (matching the
# original inlineStr cells) instead of auto-converting numeric-looking
# strings (e.g. "591.41") into real numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell 2 4 "61.193.03"
Set-TextCell 2 5 "  +0.47%  "

# Row 3
Set-TextCell 3 4 "2.924.12"
Set-TextCell 3 5 "  +0.17%  "

# Row 4
Set-TextCell 4 4 "0.999"
Set-TextCell 4 5 "  +0.10%  "

# Row 5
Set-TextCell 5 4 "591.41"
Set-TextCell 5 5 "  +0.79%  "

# Row 6
Set-TextCell 6 4 "145.94"
Set-TextCell 6 5 "  -0.96%  "

# Row 7
Set-TextCell 7 5 "  +0.05%  "

# Row 8
Set-TextCell 8 5 "  +1.08%  "

# Row 9
Set-TextCell 9 4 "2.923.67"
Set-TextCell 9 5 "  +0.00%  "

# Row 10
Set-TextCell 10 5 "  +0.99%  "

# Row 11
Set-TextCell 11 4 "0.145"
Set-TextCell 11 5 "  +0.00%  "

# Row 12
Set-TextCell 12 5 "  -1.10%  "

# Row 13
Set-TextCell 13 5 "  +1.44%  "

# Row 14
Set-TextCell 14 4 "33.86"
Set-TextCell 14 5 "  -1.64%  "

# Row 15
Set-TextCell 15 5 "  +0.00%  "

# Row 16
Set-TextCell 16 4 "3.406.27"
Set-TextCell 16 5 "  +0.13%  "

# Row 17
Set-TextCell 17 4 "61.101.35"
Set-TextCell 17 5 "  +0.43%  "

# Row 18
Set-TextCell 18 4 "6.74"
Set-TextCell 18 5 "  -1.53%  "

# Row 19
Set-TextCell 19 4 "2.920.53"
Set-TextCell 19 5 "  -0.13%  "

# Row 20
Set-TextCell 20 4 "431.65"
Set-TextCell 20 5 "  +0.78%  "

# Row 21
Set-TextCell 21 5 "  -2.07%  "

# Row 22
Set-TextCell 22 5 "  +1.77%  "

# Row 23
Set-TextCell 23 5 "  -1.59%  "

# Row 24
Set-TextCell 24 4 "80.84"
Set-TextCell 24 5 "  +0.03%  "

# Row 25
Set-TextCell 25 2 "Fetch.AI"
Set-TextCell 25 3 "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell 25 4 "2.23"
Set-TextCell 25 5 "  +2.73%  "

# Row 26
Set-TextCell 26 2 "RenderToken"
Set-TextCell 26 3 "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell 26 4 "10.79"
Set-TextCell 26 5 "  -1.65%  "

# Row 27
Set-TextCell 27 4 "12.15"
Set-TextCell 27 5 "  +2.46%  "

# Row 29
Set-TextCell 29 5 "  +8.08%  "

# Row 30
Set-TextCell 30 5 "  +0.05%  "

# Row 31
Set-TextCell 31 2 "NEARProtocol"
Set-TextCell 31 3 "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell 31 4 "7.21"
Set-TextCell 31 5 "  -1.21%  "

# Row 32
Set-TextCell 32 2 "PancakeSwap"
Set-TextCell 32 3 "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell 32 4 "2.63"
Set-TextCell 32 5 "  -0.06%  "

# Row 33
Set-TextCell 33 4 "26.61"
Set-TextCell 33 5 "  -0.89%  "

# Row 34
Set-TextCell 34 5 "  +1.51%  "

# Row 35
Set-TextCell 35 4 "0.0₃0868"
Set-TextCell 35 5 "  +3.21%  "

# Row 36
Set-TextCell 36 4 "1.01"
Set-TextCell 36 5 "  +0.72%  "

# Row 37
Set-TextCell 37 4 "3.12"
Set-TextCell 37 5 "  +4.25%  "

# Row 38
Set-TextCell 38 5 "  -0.88%  "

# Row 39
Set-TextCell 39 4 "49.95"
Set-TextCell 39 5 "  +0.24%  "

# Row 40
Set-TextCell 40 5 "  -1.05%  "

# Row 41
Set-TextCell 41 5 "  -0.66%  "

# Row 42
Set-TextCell 42 4 "8.62"
Set-TextCell 42 5 "  -1.69%  "

# Row 43
Set-TextCell 43 5 "  -0.23%  "

# Row 44
Set-TextCell 44 4 "39.66"
Set-TextCell 44 5 "  -4.93%  "

# Row 45
Set-TextCell 45 4 "381.16"
Set-TextCell 45 5 "  +1.96%  "

# Row 46
Set-TextCell 46 4 "0.0349"
Set-TextCell 46 5 "  +0.87%  "

# Row 47
Set-TextCell 47 4 "2.708.44"
Set-TextCell 47 5 "  +1.88%  "

# Row 48
Set-TextCell 48 4 "130.91"
Set-TextCell 48 5 "  -1.38%  "

# Row 49
Set-TextCell 49 5 "  +0.01%  "

# Row 50
Set-TextCell 50 4 "24.42"
Set-TextCell 50 5 "  -3.55%  "

# Row 51
Set-TextCell 51 5 "  +0.39%  "
